{"js": "const replacements = [\n  [\"2024-07-30 Tuesday\", \"2024-07-31 Wednesday\"],\n  [\"179\u00f75=35, 4\", \"778\u00f76=129, 4\"],\n  [\"137\u00f72=68, 1\", \"576\u00f72=288, 0\"],\n  [\"856\u00f79=95, 1\", \"309\u00f75=61, 4\"],\n  [\"161\u00f78=20, 1\", \"825\u00f78=103, 1\"],\n  [\"519\u00f78=64, 7\", \"265\u00f75=53, 0\"],\n  [\"106\u00f72=53, 0\", \"913\u00f77=130, 3\"],\n  [\"978\u00f72=489, 0\", \"300\u00f75=60, 0\"],\n  [\"187\u00f72=93, 1\", \"529\u00f76=88, 1\"],\n  [\"487\u00f74=121, 3\", \"132\u00f76=22, 0\"],\n  [\"545\u00f75=109, 0\", \"573\u00f75=114, 3\"],\n  [\"531\u00f72=265, 1\", \"207\u00f77=29, 4\"],\n  [\"558\u00f79=62, 0\", \"485\u00f72=242, 1\"],\n  [\"510\u00f73=170, 0\", \"402\u00f72=201, 0\"],\n  [\"414\u00f79=46, 0\", \"844\u00f77=120, 4\"],\n  [\"439\u00f74=109, 3\", \"429\u00f75=85, 4\"],\n  [\"421\u00f73=140, 1\", \"786\u00f76=131, 0\"],\n  [\"331\u00f79=36, 7\", \"234\u00f74=58, 2\"],\n  [\"154\u00f79=17, 1\", \"176\u00f74=44, 0\"],\n  [\"154\u00f75=30, 4\", \"389\u00f78=48, 5\"],\n  [\"645\u00f76=107, 3\", \"630\u00f78=78, 6\"],\n  [\"495\u00f78=61, 7\", \"469\u00f79=52, 1\"],\n  [\"916\u00f74=229, 0\", \"411\u00f73=137, 0\"],\n  [\"153\u00f78=19, 1\", \"638\u00f77=91, 1\"],\n  [\"461\u00f79=51, 2\", \"200\u00f77=28, 4\"],\n  [\"876\u00f75=175, 1\", \"830\u00f72=415, 0\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-30 Tuesday\", \"2024-07-31 Wednesday\"),\n    @(\"179\u00f75=35, 4\", \"778\u00f76=129, 4\"),\n    @(\"137\u00f72=68, 1\", \"576\u00f72=288, 0\"),\n    @(\"856\u00f79=95, 1\", \"309\u00f75=61, 4\"),\n    @(\"161\u00f78=20, 1\", \"825\u00f78=103, 1\"),\n    @(\"519\u00f78=64, 7\", \"265\u00f75=53, 0\"),\n    @(\"106\u00f72=53, 0\", \"913\u00f77=130, 3\"),\n    @(\"978\u00f72=489, 0\", \"300\u00f75=60, 0\"),\n    @(\"187\u00f72=93, 1\", \"529\u00f76=88, 1\"),\n    @(\"487\u00f74=121, 3\", \"132\u00f76=22, 0\"),\n    @(\"545\u00f75=109, 0\", \"573\u00f75=114, 3\"),\n    @(\"531\u00f72=265, 1\", \"207\u00f77=29, 4\"),\n    @(\"558\u00f79=62, 0\", \"485\u00f72=242, 1\"),\n    @(\"510\u00f73=170, 0\", \"402\u00f72=201, 0\"),\n    @(\"414\u00f79=46, 0\", \"844\u00f77=120, 4\"),\n    @(\"439\u00f74=109, 3\", \"429\u00f75=85, 4\"),\n    @(\"421\u00f73=140, 1\", \"786\u00f76=131, 0\"),\n    @(\"331\u00f79=36, 7\", \"234\u00f74=58, 2\"),\n    @(\"154\u00f79=17, 1\", \"176\u00f74=44, 0\"),\n    @(\"154\u00f75=30, 4\", \"389\u00f78=48, 5\"),\n    @(\"645\u00f76=107, 3\", \"630\u00f78=78, 6\"),\n    @(\"495\u00f78=61, 7\", \"469\u00f79=52, 1\"),\n    @(\"916\u00f74=229, 0\", \"411\u00f73=137, 0\"),\n    @(\"153\u00f78=19, 1\", \"638\u00f77=91, 1\"),\n    @(\"461\u00f79=51, 2\", \"200\u00f77=28, 4\"),\n    @(\"876\u00f75=175, 1\", \"830\u00f72=415, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute(\n        $findText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $replaceText,\n        2\n    )\n}\n\nWrite-Output \"done\""}
